$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GA")
$ws.Activate()

# Update existing C column totals
$ws.Range("C2").Value = 6952
$ws.Range("C3").Value = 6953

# Add new data in E/F columns and formula in G
$ws.Range("E2").Value = 2001
$ws.Range("F2").Value = 4951
$ws.Range("F2").Borders.LineStyle = 1
$ws.Range("F2").Borders.Weight = 2
$ws.Range("G2").Formula = "=E2+F2"

$ws.Range("E3").Value = 2002
$ws.Range("F3").Value = 4951
$ws.Range("F3").Borders.LineStyle = 1
$ws.Range("F3").Borders.Weight = 2
$ws.Range("G3").Formula = "=E3+F3"

# Update selection to match target state
$ws.Range("I20").Select() | Out-Null
